# "new run with new features"
# Updates the classification_evaluation sheet with results from a new
# experiment run: the set of classifiers compared is the same, but the
# row order (classifier labels in column A) and the accuracy /
# f1 / threshold scores (columns B/C/D) reflect the new run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{Row=2;  A="SVC rbf";                     B=0.717;  C=0.717;  D=0.9272},
    @{Row=3;  A="MLP-deep";                    B=0.6813; C=0.6813; D=0.9245},
    @{Row=4;  A="MLP 128";                     B=0.6827; C=0.6827; D=0.9217},
    @{Row=5;  A="MLP 16";                      B=0.6841; C=0.6841; D=0.9217},
    @{Row=6;  A="LinearSVC";                   B=0.6731; C=0.6731; D=0.919},
    @{Row=7;  A="SVC poly";                    B=0.6772; C=0.6772; D=0.919},
    @{Row=8;  A="GradientBoostingClassifier";  B=0.7033; C=0.7033; D=0.9148},
    @{Row=9;  A="LogisticRegressionCV";        B=0.6813; C=0.6813; D=0.9148},
    @{Row=10; A="MLP 32";                      B=0.6813; C=0.6813; D=0.9121},
    @{Row=11; A="MLP 64";                      B=0.6786; C=0.6786; D=0.9121},
    @{Row=12; A="RandomForestClassifier";      B=0.6772; C=0.6772; D=0.9107},
    @{Row=13; A="SVC sigmoid";                 B=0.6016; C=0.6016; D=0.8503}
)

foreach ($row in $rows) {
    $ws.Cells.Item($row.Row, 1).Value = $row.A
    $ws.Cells.Item($row.Row, 2).Value = $row.B
    $ws.Cells.Item($row.Row, 3).Value = $row.C
    $ws.Cells.Item($row.Row, 4).Value = $row.D
}
